$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "306.76"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "0.85%"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "36.28"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "2.54%"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.098"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "0.06%"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.08105"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "1.49%"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.966"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "0.18%"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "7.748"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "-1.20%"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.9326"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "0.74%"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1478"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "35.26%"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.1922"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "1.50%"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.09100"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "-3.68%"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.03527"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "-3.95%"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.09821"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-1.38%"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001450"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "1.41%"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.005791"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "-0.18%"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "1.95%"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.103"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "-0.82%"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.933"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "1.27%"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.3427"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "0.33%"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.1299"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "-0.94%"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.045"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "-1.38%"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.2397"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "8.86%"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04526"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "-0.06%"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001211"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "-1.59%"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004893"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "4.42%"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0001245"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "-0.98%"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0004435"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "-0.88%"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01989"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "4.52%"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.04855"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "2.40%"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.01112"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "14.92%"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.007564"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "-0.35%"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.1374"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "2.06%"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.002088"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "-1.86%"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.009888"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "-12.88%"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00006211"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "-1.41%"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.00000000754"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "-0.09%"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "0.44%"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.001189"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "-9.10%"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.00002110"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "-0.09%"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0002010"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "-0.09%"
